# Apply the price/volume/coin-identity updates from the Fri Mar 8 05:37:32 UTC 2024
# GitHub Actions "Updated cryptos list" commit.
#
# All data cells on the sheet are plain text ("inlineStr") cells, including the
# Price column (D) which holds numeric-looking strings (e.g. "66.858.76",
# "1.00", "0.140"). Excel auto-converts a bare numeric-looking Value assignment
# into a real number, which would silently change the cell type and drop
# formatting such as trailing zeros. To keep these as text we assign with a
# leading apostrophe (Excel's literal-text marker) and then reset the cell
# style back to "Normal" so no stray quote-prefix style lingers on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "'66.858.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.43%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "'3.881.14"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.70%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.16%  "

# Row 5: BNB
$ws.Range("D5").Value = "'468.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +9.63%  "

# Row 6: Solana
$ws.Range("D6").Value = "'143.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.84%  "

# Row 7: XRP
$ws.Range("D7").Value = "'0.622"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.48%  "

# Row 8: USDC
$ws.Range("D8").Value = "'0.998"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.19%  "

# Row 9: Cardano
$ws.Range("D9").Value = "'0.731"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.29%  "

# Row 10: Dogecoin
$ws.Range("E10").Value = "  +7.83%  "

# Row 11: ShibaInu
$ws.Range("D11").Value = "'0.0000329"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +9.42%  "

# Row 12: Avalanche
$ws.Range("D12").Value = "'42.77"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.39%  "

# Row 13: WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "'4.507.95"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.53%  "

# Row 14: Polkadot
$ws.Range("E14").Value = "  -0.99%  "

# Row 15: Uniswap
$ws.Range("D15").Value = "'14.89"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.25%  "

# Row 16: WrappedEther
$ws.Range("D16").Value = "'3.907.04"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.81%  "

# Row 17: TRON
$ws.Range("E17").Value = "  -0.36%  "

# Row 18: Chainlink
$ws.Range("D18").Value = "'19.77"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.19%  "

# Row 19: Polygon
$ws.Range("E19").Value = "  +3.65%  "

# Row 20: WrappedBTC
$ws.Range("D20").Value = "'67.046.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.36%  "

# Row 21: BitcoinCash
$ws.Range("D21").Value = "'427.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.12%  "

# Row 22: InternetComputer(DFINITY)
$ws.Range("D22").Value = "'14.63"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.22%  "

# Row 23: ImmutableX
$ws.Range("E23").Value = "  +3.59%  "

# Row 24: Litecoin
$ws.Range("D24").Value = "'87.83"
$ws.Range("D24").Style = "Normal"

# Row 25: PancakeSwap
$ws.Range("D25").Value = "'3.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.98%  "

# Row 26: EthereumClassic
$ws.Range("D26").Value = "'38.19"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.21%  "

# Row 27: LEO
$ws.Range("D27").Value = "'5.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.23%  "

# Row 28: Filecoin
$ws.Range("D28").Value = "'9.94"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.31%  "

# Row 29: RenderToken
$ws.Range("D29").Value = "'9.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.18%  "

# Row 30: Bittensor
$ws.Range("D30").Value = "'728.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.08%  "

# Row 31: Cosmos
$ws.Range("D31").Value = "'13.69"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.15%  "

# Row 32: Hedera
$ws.Range("E32").Value = "  -0.45%  "

# Row 33: Toncoin
$ws.Range("E33").Value = "  +0.43%  "

# Row 34: InjectiveProtocol
$ws.Range("D34").Value = "'43.16"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.55%  "

# Row 35: Kaspa
$ws.Range("E35").Value = "  +5.65%  "

# Row 36: OKB
$ws.Range("D36").Value = "'57.46"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.65%  "

# Row 37: Dai
$ws.Range("E37").Value = "  -0.14%  "

# Row 38: NEARProtocol
$ws.Range("D38").Value = "'5.37"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.36%  "

# Row 39: PEPE
$ws.Range("D39").Value = "'0.0₃0765"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +14.69%  "

# Row 40: VeChain
$ws.Range("D40").Value = "'0.0473"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.99%  "

# Row 41: ThetaToken
$ws.Range("D41").Value = "'3.08"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.40%  "

# Row 42: Stellar -> FirstDigitalUSD
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "'1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.06%  "

# Row 43: FirstDigitalUSD -> Stellar
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").Value = "'0.140"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.30%  "

# Row 44: TheGraph
$ws.Range("D44").Value = "'0.335"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.35%  "

# Row 45: WEMIXToken -> Fetch.AI
$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").Value = "'2.54"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.94%  "

# Row 46: Fetch.AI -> WEMIXToken
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "'2.79"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.56%  "

# Row 47: ARBITRUM
$ws.Range("E47").Value = "  +5.01%  "

# Row 48: LidoDAOToken
$ws.Range("D48").Value = "'3.37"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.96%  "

# Row 49: ApeXProtocol
$ws.Range("D49").Value = "'3.11"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.17%  "

# Row 50: Monero
$ws.Range("D50").Value = "'143.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.73%  "

# Row 51: Stacks
$ws.Range("E51").Value = "  +3.31%  "
